$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.618.81"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.821.61"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "305.64"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "0.4680"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("D8").Value = "0.3592"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "46.21"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.07135"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.9026"
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07802"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "19.41"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.813.60"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.260"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "6.336"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "87.36"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000008555"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.008"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.677.60"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "14.19"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.014"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.55"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.936"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "151.85"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "17.90"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "1.972"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "113.57"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "4.809"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08807"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "3.144"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "2.775"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7305"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.440"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.124"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.079"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01929"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.920"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05112"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5058"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.832"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1499"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "8.000"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4677"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.03"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "99.03"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.561"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06003"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "63.85"
$ws.Range("E51").Value = "  -0.53%  "
